$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column O: "Addressing mode(Optional)" with Static/Dynamic values ---
# (write in this order so the new shared-string table is built as
#  32=Addressing mode(Optional), 33=Static, 34=Dynamic - matching the target workbook)
$ws.Range("O1").Value = "Addressing mode(Optional)"
$ws.Range("O3").Value = "Static"
$ws.Range("O2").Value = "Dynamic"

# Set the new column's width
$ws.Columns.Item(15).ColumnWidth = 26.4

# --- Update the "Speed limit M/s(Optional)" header cell (N1) ---
# Split the zero-width-space run out into its own run using the MS Gothic font,
# matching the new rich-text layout: "Speed " + (zero width spaces in MS Gothic) + "limit M/s(Optional)"
$speedCell = $ws.Range("N1")
$zwsp = $speedCell.Characters(7, 2)
$zwsp.Font.Name = "MS Gothic"
$zwsp.Font.Size = 12

# --- Update sheet view: clear the scrolled viewport and move the selection ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("G17").Select()
